# Update the cryptocurrency price table (rows 2-51) to reflect the
# latest scrape: prices, 1h volume %, and (where the ranking shifted)
# coin name / link columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.670.54"
$ws.Range("E2").Value = "  -5.95%  "

# Row 3
$ws.Range("D3").Value = "3.261.82"
$ws.Range("E3").Value = "  -7.53%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").Value = "'517.80"
$ws.Range("E5").Value = "  -6.99%  "

# Row 6
$ws.Range("D6").Value = "'172.31"
$ws.Range("E6").Value = "  -15.31%  "

# Row 7
$ws.Range("D7").Value = "'0.599"
$ws.Range("E7").Value = "  -1.74%  "

# Row 8
$ws.Range("D8").Value = "3.256.26"
$ws.Range("E8").Value = "  -7.34%  "

# Row 9
$ws.Range("E9").Value = "  +0.13%  "

# Row 10
$ws.Range("D10").Value = "'0.601"
$ws.Range("E10").Value = "  -8.71%  "

# Row 11
$ws.Range("D11").Value = "'56.03"
$ws.Range("E11").Value = "  -13.23%  "

# Row 12
$ws.Range("D12").Value = "'0.131"
$ws.Range("E12").Value = "  -8.94%  "

# Row 13
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -6.45%  "

# Row 14
$ws.Range("D14").Value = "'8.94"
$ws.Range("E14").Value = "  -9.93%  "

# Row 15
$ws.Range("D15").Value = "3.805.73"
$ws.Range("E15").Value = "  -6.69%  "

# Row 16
$ws.Range("D16").Value = "3.273.77"
$ws.Range("E16").Value = "  -6.93%  "

# Row 17
$ws.Range("D17").Value = "'0.115"
$ws.Range("E17").Value = "  -7.02%  "

# Row 18
$ws.Range("D18").Value = "63.677.51"
$ws.Range("E18").Value = "  -5.54%  "

# Row 19
$ws.Range("D19").Value = "'17.17"
$ws.Range("E19").Value = "  -7.53%  "

# Row 20
$ws.Range("D20").Value = "'10.93"
$ws.Range("E20").Value = "  -7.96%  "

# Row 21
$ws.Range("D21").Value = "'0.943"
$ws.Range("E21").Value = "  -8.62%  "

# Row 22
$ws.Range("D22").Value = "'369.36"
$ws.Range("E22").Value = "  -6.23%  "

# Row 23
$ws.Range("D23").Value = "'3.72"
$ws.Range("E23").Value = "  -7.60%  "

# Row 24
$ws.Range("D24").Value = "'79.50"
$ws.Range("E24").Value = "  -4.70%  "

# Row 25
$ws.Range("D25").Value = "'10.84"
$ws.Range("E25").Value = "  -11.13%  "

# Row 26
$ws.Range("D26").Value = "'3.83"
$ws.Range("E26").Value = "  -2.37%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.62"
$ws.Range("E27").Value = "  -7.85%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'11.20"
$ws.Range("E28").Value = "  -8.64%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'8.18"
$ws.Range("E29").Value = "  -8.19%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'28.39"
$ws.Range("E30").Value = "  -8.93%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'635.51"
$ws.Range("E31").Value = "  -11.18%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.53"
$ws.Range("E32").Value = "  -8.39%  "

# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.10"
$ws.Range("E33").Value = "  -5.83%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'58.40"
$ws.Range("E34").Value = "  -8.88%  "

# Row 35
$ws.Range("D35").Value = "'0.104"
$ws.Range("E35").Value = "  -7.34%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.15%  "

# Row 37
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'36.08"
$ws.Range("E37").Value = "  -7.01%  "

# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.380"
$ws.Range("E38").Value = "  -5.02%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.50%  "

# Row 40
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0686"
$ws.Range("E40").Value = "  -0.62%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.895.76"
$ws.Range("E41").Value = "  -5.47%  "

# Row 42
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "  -7.10%  "

# Row 43
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.41"
$ws.Range("E43").Value = "  -7.29%  "

# Row 44
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'2.64"
$ws.Range("E44").Value = "  -12.45%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.63"
$ws.Range("E45").Value = "  -4.58%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0390"
$ws.Range("E46").Value = "  -4.53%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'2.99"
$ws.Range("E47").Value = "  +4.37%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.124"
$ws.Range("E48").Value = "  -2.76%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.72"
$ws.Range("E49").Value = "  +3.04%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'134.00"
$ws.Range("E50").Value = "  -2.98%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'2.33"
$ws.Range("E51").Value = "  -15.58%  "
